$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header column (H), matching the existing header formatting
# (bold, bordered, centered) by copying the adjacent header cell's format.
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Per-row label values (0/1) for the new column.
$labels = @(0, 0, 1, 1, 1, 1, 1, 0, 0, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
